$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D values entered first (matches shared-string insertion order)
$ws.Range("D2").Value = "teste"
$ws.Range("D3").Value = "alo"
$ws.Range("D4").Value = "figuras"

# Update header row (columns stay in place, text changes) then add D1 header
$ws.Range("A1").Value = "Cod"
$ws.Range("B1").Value = "Disciplina"
$ws.Range("C1").Value = "Carga Horaria"
$ws.Range("D1").Value = "Opcionais"

# Match the column D width seen in the target workbook (closest value the
# engine's pixel-quantized ColumnWidth model can reproduce for 12.140625)
$ws.Columns.Item(4).ColumnWidth = 11.333333333333332

# Update the active selection to D1, matching the saved view state
$ws.Range("D1").Select()
